# Final Push Part 1
# Insert a new row for the "Not Taken" / "TP 500 Team Project Full" option
# between the "Prerequisites" row (9) and "Aim of module" row (previously 10,
# now shifted to 11), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10; existing rows 10-25 shift to 11-26.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row.
$ws.Range("A10").Value = "Not Taken"
$ws.Range("B10").Value = "TP 500 Team Project Full"

# Move the active selection to the newly edited cell, matching the saved
# workbook state.
$ws.Range("B10").Select()
